$d = $word.ActiveDocument

# The paragraph originally reads "...她简直好不拒绝，甚至..." — fix the
# typo "好不" -> "毫不" ("她简直毫不拒绝" = "she didn't hesitate at all").
$d.Content.Find.Execute("好不拒绝", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "毫不拒绝", 2)

# Locate the point right after "毫不" (i.e. right before "拒绝") — this is
# where the run needs to be split in two.
$r = $d.Content
$r.Find.Execute("毫不", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "", 0)
$splitPoint = $r.End

# The "_GoBack" bookmark currently sits at the end of the paragraph's run.
# Move it to the split point: deleting it and re-adding it there forces the
# run to be split into two runs (one ending in "...毫不", one starting with
# "拒绝...") with the bookmark sitting between them, matching the edit that
# was made when the author last typed at that spot.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()
$splitRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $splitRange)
